$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking strings (e.g. "27.528.44", "1.00")
# that must remain plain text, matching the source inlineStr cells. Forcing
# a text number format before assignment (and restoring the default style
# afterwards) prevents Excel from auto-converting them into real numbers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '27.528.44'
$ws.Range('E2').Value = '  -0.08%  '
Set-TextValue $ws.Range('D3') '1.618.45'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue $ws.Range('D5') '210.92'
$ws.Range('E5').Value = '  -0.70%  '
Set-TextValue $ws.Range('D6') '0.527'
$ws.Range('E6').Value = '  -1.80%  '
$ws.Range('E7').Value = '  -0.03%  '
Set-TextValue $ws.Range('D8') '22.64'
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('E9').Value = '  +2.47%  '
$ws.Range('E10').Value = '  +0.20%  '
Set-TextValue $ws.Range('D11') '0.0887'
$ws.Range('E11').Value = '  -0.26%  '
Set-TextValue $ws.Range('D12') '1.846.35'
$ws.Range('E12').Value = '  -1.42%  '
Set-TextValue $ws.Range('D13') '1.616.80'
$ws.Range('E13').Value = '  -1.55%  '
$ws.Range('E14').Value = '  -0.07%  '
Set-TextValue $ws.Range('D15') '0.551'
$ws.Range('E15').Value = '  -1.83%  '
Set-TextValue $ws.Range('D16') '64.95'
$ws.Range('E16').Value = '  +1.57%  '
Set-TextValue $ws.Range('D17') '27.527.38'
$ws.Range('E17').Value = '  +0.02%  '
Set-TextValue $ws.Range('D18') '230.10'
$ws.Range('E18').Value = '  +0.96%  '
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('E20').Value = '  -0.82%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('E22').Value = '  +0.05%  '
Set-TextValue $ws.Range('D23') '10.10'
$ws.Range('E23').Value = '  +1.32%  '
$ws.Range('E24').Value = '  +7.32%  '
Set-TextValue $ws.Range('D25') '149.39'
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('E26').Value = '  -0.93%  '
$ws.Range('B27').Value = 'BinanceUSD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range('D27') '1.00'
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D28') '6.80'
$ws.Range('E28').Value = '  -2.01%  '
Set-TextValue $ws.Range('D29') '15.54'
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('E31').Value = '  -0.58%  '
$ws.Range('E32').Value = '  -0.85%  '
Set-TextValue $ws.Range('D33') '1.442.57'
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('E34').Value = '  -3.27%  '
$ws.Range('E35').Value = '  -3.02%  '
$ws.Range('E36').Value = '  -0.31%  '
Set-TextValue $ws.Range('D37') '0.934'
$ws.Range('E37').Value = '  +3.02%  '
$ws.Range('E38').Value = '  -2.03%  '
$ws.Range('E39').Value = '  +0.33%  '
$ws.Range('E40').Value = '  -1.64%  '
Set-TextValue $ws.Range('D41') '69.14'
$ws.Range('E41').Value = '  +6.40%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('E43').Value = '  -3.05%  '
$ws.Range('E44').Value = '  -0.34%  '
Set-TextValue $ws.Range('D45') '5.40'
$ws.Range('E45').Value = '  -1.80%  '
$ws.Range('E46').Value = '  -1.92%  '
Set-TextValue $ws.Range('D47') '1.758.89'
$ws.Range('E47').Value = '  -1.35%  '
Set-TextValue $ws.Range('D48') '1.69'
$ws.Range('E48').Value = '  +0.44%  '
Set-TextValue $ws.Range('D49') '86.46'
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('E50').Value = '  -1.74%  '
Set-TextValue $ws.Range('D51') '0.0996'
$ws.Range('E51').Value = '  +1.45%  '
